# The "Recorded By" column (G) lists the users/accounts that recorded a
# given attendance session, as a comma-separated string. Previously the
# literal token "System" was always placed first. This edit normalizes
# the ordering so "System" (when present as the first entry) is moved to
# the end of the list instead, e.g.:
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $value = $cell.Value()

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = New-Object System.Collections.Generic.List[string]
    foreach ($p in $value.Split(",")) {
        [void]$parts.Add($p.Trim())
    }

    if ($parts.Count -gt 1 -and $parts[0].Equals("System")) {
        $first = $parts[0]
        [void]$parts.RemoveAt(0)
        [void]$parts.Add($first)
        $cell.Value = [string]::Join(", ", $parts)
    }
}
